# Insert a new data row for 2026/01/14 (time slot 17) right before the
# 2026/12/29 block, which pushes the existing rows 637..678 down to 638..679.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 637:678 down to 638:679 by inserting a new row at 637.
$ws.Rows("637:637").Insert()

# Populate the newly inserted row with the new data point. The date column
# is stored as literal text (not a real date) elsewhere in the sheet, so
# force text formatting before assigning to avoid Excel auto-converting the
# "2026/01/14" string into a date serial number; restore the default
# (Normal) style afterwards so the new row matches the surrounding rows.
$ws.Cells.Item(637, 1).NumberFormat = "@"
$ws.Cells.Item(637, 1).Value = "2026/01/14"
$ws.Cells.Item(637, 1).Style = "Normal"
$ws.Cells.Item(637, 2).Value = "水"
$ws.Cells.Item(637, 3).Value = 17
$ws.Cells.Item(637, 4).Value = 201
